$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 7-16 first (from the bottom up) so row indices of the
# remaining rows (1-6) are unaffected while we still have 16 rows.
$ws.Rows("7:16").Delete()

# Update rows 2-6 with the new combined tuple-style strings.
$ws.Range("A2").Value = "('Forest', ['Basic Land " + [char]0x2014 + " Forest', '({T}: Add {G}.)'])"
$ws.Range("A3").Value = "('Island', ['Basic Land " + [char]0x2014 + " Island', '({T}: Add {U}.)'])"
$ws.Range("A4").Value = "('Mountain', ['Basic Land " + [char]0x2014 + " Mountain', '({T}: Add {R}.)'])"
$ws.Range("A5").Value = "('Plains', ['Basic Land " + [char]0x2014 + " Plains', '({T}: Add {W}.)'])"
$ws.Range("A6").Value = "('Swamp', ['Basic Land " + [char]0x2014 + " Swamp', '({T}: Add {B}.)'])"
